# OrderInvoiceTemplate.xlsx fix:
#  - Insert a new "Price Discount" column between "Discount in %" and "VAT"
#  - Add a new data row describing a "Home Elevator" order
#  - Resize the new columns and update the active selection

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column D ("Price Discount"). This shifts the existing
# "VAT" column (old D) to E and "Price with VAT" (old E) to F.
$ws.Columns("D").Insert()
$ws.Range("D1").Value = "Price Discount"

# New order-data row
$ws.Range("A2").Value = "Home Elevator"
$ws.Range("B2").Value = 23440
$ws.Range("C2").Value = 0
$ws.Range("D2").Value = 0
$ws.Range("E2").Value = 20
$ws.Range("F2").Value = 0

# Column widths for the new "Price Discount" (D) and "Price with VAT" (F) columns
$ws.Columns("D").ColumnWidth = 16.8
$ws.Columns("F").ColumnWidth = 15.1

# Restore the selection to match the saved view state
$ws.Range("G4").Select()
